$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 676, pushing the existing records (rows 676-747)
# down to rows 677-748.
$ws.Rows("676:676").Insert()

# Populate the newly inserted row 676 with the new weekly record.
$ws.Range("A676").Value = 10
$ws.Range("B676").Value = "Vega Modelo de Temuco"
$ws.Range("C676").Value = "La Araucanía"
$ws.Range("D676").Value = 45194
$ws.Range("E676").Value = 9
$ws.Range("F676").Value = 100112028
$ws.Range("G676").Value = "Sandia"
$ws.Range("H676").Value = "Sin especificar"
$ws.Range("I676").Value = "Primera"
$ws.Range("J676").Value = 800
$ws.Range("K676").Value = 1000
$ws.Range("L676").Value = 1000
$ws.Range("M676").Value = 1000
$ws.Range("N676").Value = "$/kilo (volumen en unidades)"
$ws.Range("O676").Value = "Perú"
$ws.Range("P676").Value = 1000
$ws.Range("Q676").Value = 1
$ws.Range("R676").Value = "Hortaliza"
